# Updated symbol list on Sun Jan  8 04:12:36 UTC 2023 with GitHub Actions
# Refreshes Price / Volume(1h) / Hora columns for the cryptos sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "261.10" },
    @{ Cell = "E2"; Value = "-0.29%" },
    @{ Cell = "G2"; Value = "4" },
    @{ Cell = "E3"; Value = "-1.02%" },
    @{ Cell = "G3"; Value = "4" },
    @{ Cell = "D4"; Value = "4.701" },
    @{ Cell = "E4"; Value = "-0.82%" },
    @{ Cell = "G4"; Value = "4" },
    @{ Cell = "D5"; Value = "0.06226" },
    @{ Cell = "E5"; Value = "2.65%" },
    @{ Cell = "G5"; Value = "4" },
    @{ Cell = "D6"; Value = "6.739" },
    @{ Cell = "E6"; Value = "0.37%" },
    @{ Cell = "G6"; Value = "4" },
    @{ Cell = "D7"; Value = "0.8514" },
    @{ Cell = "E7"; Value = "-1.78%" },
    @{ Cell = "G7"; Value = "4" },
    @{ Cell = "D8"; Value = "0.9103" },
    @{ Cell = "E8"; Value = "-1.41%" },
    @{ Cell = "G8"; Value = "4" },
    @{ Cell = "D9"; Value = "0.1397" },
    @{ Cell = "E9"; Value = "-0.87%" },
    @{ Cell = "G9"; Value = "4" },
    @{ Cell = "D10"; Value = "0.04855" },
    @{ Cell = "E10"; Value = "-2.34%" },
    @{ Cell = "G10"; Value = "4" },
    @{ Cell = "D11"; Value = "0.07090" },
    @{ Cell = "E11"; Value = "-1.01%" },
    @{ Cell = "G11"; Value = "4" },
    @{ Cell = "D12"; Value = "0.03136" },
    @{ Cell = "E12"; Value = "2.97%" },
    @{ Cell = "G12"; Value = "4" },
    @{ Cell = "D13"; Value = "0.09066" },
    @{ Cell = "E13"; Value = "-0.58%" },
    @{ Cell = "G13"; Value = "4" },
    @{ Cell = "E14"; Value = "0.29%" },
    @{ Cell = "G14"; Value = "4" },
    @{ Cell = "D15"; Value = "0.0006175" },
    @{ Cell = "E15"; Value = "1.32%" },
    @{ Cell = "G15"; Value = "4" },
    @{ Cell = "D16"; Value = "0.006045" },
    @{ Cell = "E16"; Value = "-2.40%" },
    @{ Cell = "G16"; Value = "4" },
    @{ Cell = "D17"; Value = "3.450" },
    @{ Cell = "E17"; Value = "0.05%" },
    @{ Cell = "G17"; Value = "4" },
    @{ Cell = "D18"; Value = "3.175" },
    @{ Cell = "E18"; Value = "0.17%" },
    @{ Cell = "G18"; Value = "4" },
    @{ Cell = "D19"; Value = "2.166" },
    @{ Cell = "E19"; Value = "-0.49%" },
    @{ Cell = "G19"; Value = "4" },
    @{ Cell = "G20"; Value = "4" },
    @{ Cell = "D21"; Value = "0.1310" },
    @{ Cell = "E21"; Value = "0.96%" },
    @{ Cell = "G21"; Value = "4" },
    @{ Cell = "D22"; Value = "4.128" },
    @{ Cell = "E22"; Value = "0.94%" },
    @{ Cell = "G22"; Value = "4" },
    @{ Cell = "D23"; Value = "0.04247" },
    @{ Cell = "E23"; Value = "-0.32%" },
    @{ Cell = "G23"; Value = "4" },
    @{ Cell = "D24"; Value = "0.001216" },
    @{ Cell = "E24"; Value = "-0.20%" },
    @{ Cell = "G24"; Value = "4" },
    @{ Cell = "D25"; Value = "0.004085" },
    @{ Cell = "E25"; Value = "4.40%" },
    @{ Cell = "G25"; Value = "4" },
    @{ Cell = "E26"; Value = "0.02%" },
    @{ Cell = "G26"; Value = "4" },
    @{ Cell = "D27"; Value = "0.0001640" },
    @{ Cell = "E27"; Value = "4.36%" },
    @{ Cell = "G27"; Value = "4" },
    @{ Cell = "G28"; Value = "4" },
    @{ Cell = "G29"; Value = "4" },
    @{ Cell = "G30"; Value = "4" },
    @{ Cell = "G31"; Value = "4" },
    @{ Cell = "G32"; Value = "4" },
    @{ Cell = "G33"; Value = "4" },
    @{ Cell = "G34"; Value = "4" },
    @{ Cell = "G35"; Value = "4" },
    @{ Cell = "G36"; Value = "4" },
    @{ Cell = "G37"; Value = "4" },
    @{ Cell = "G38"; Value = "4" },
    @{ Cell = "G39"; Value = "4" },
    @{ Cell = "D40"; Value = "0.03937" },
    @{ Cell = "E40"; Value = "1.44%" },
    @{ Cell = "G40"; Value = "4" },
    @{ Cell = "D41"; Value = "0.1112" },
    @{ Cell = "E41"; Value = "-0.25%" },
    @{ Cell = "G41"; Value = "4" },
    @{ Cell = "D42"; Value = "0.004127" },
    @{ Cell = "E42"; Value = "-0.15%" },
    @{ Cell = "G42"; Value = "4" },
    @{ Cell = "D43"; Value = "0.002211" },
    @{ Cell = "E43"; Value = "0.29%" },
    @{ Cell = "G43"; Value = "4" },
    @{ Cell = "D44"; Value = "0.01388" },
    @{ Cell = "E44"; Value = "-6.98%" },
    @{ Cell = "G44"; Value = "4" },
    @{ Cell = "D45"; Value = "0.00005122" },
    @{ Cell = "E45"; Value = "-4.25%" },
    @{ Cell = "G45"; Value = "4" },
    @{ Cell = "E46"; Value = "0.01%" },
    @{ Cell = "G46"; Value = "4" },
    @{ Cell = "D47"; Value = "0.03401" },
    @{ Cell = "E47"; Value = "-37.67%" },
    @{ Cell = "G47"; Value = "4" },
    @{ Cell = "D48"; Value = "0.06764" },
    @{ Cell = "E48"; Value = "-50.01%" },
    @{ Cell = "G48"; Value = "4" },
    @{ Cell = "D49"; Value = "0.00002101" },
    @{ Cell = "E49"; Value = "0.01%" },
    @{ Cell = "G49"; Value = "4" },
    @{ Cell = "D50"; Value = "0.0002001" },
    @{ Cell = "E50"; Value = "0.01%" },
    @{ Cell = "G50"; Value = "4" },
    @{ Cell = "G51"; Value = "4" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Force text storage so numeric-looking / percentage-looking strings
    # (e.g. "261.10", "-0.29%", "4") are kept as literal text, matching
    # the sheet's existing inline-string column format instead of being
    # auto-converted into numbers/percentages by Excel's type inference.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}
